$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 41: "wrap Box" packing product becomes "Butter Papper" ---
$ws.Range("B41").Value = "Packing_Product_5_ Butter Papper"
$ws.Range("C41").Value = "Butter Papper"
$ws.Range("D41").Value = "Central market"
$ws.Range("I41").Value = "Butter Papper.jpg"

# --- Add new row 42: "Polythene bag" packing product ---
$ws.Range("B42").Value = "Packing_Product_6_ Polythene bag"
$ws.Range("I42").Value = "Polythene bag.jpg"
$ws.Range("C42").Value = "Polythene bag"
$ws.Range("H42").Value = "Packing Product_6"
$ws.Range("D42").Value = "Central market"
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = "Packing"
$ws.Range("G42").Value = "Packing"

# --- Update view state to match the saved workbook ---
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("G31").Select() | Out-Null
